$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: ParticipantsTab query - updated to new participant query
$ws.Range("B2").Value = 'MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in [''XLSX'']
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY p.participant_id'

# B3: SamplesTab query - updated to new sample query
$ws.Range("B3").Value = 'MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE f.file_type in [''XLSX'']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(p.participant_id,'''') as `Participant ID`,
    coalesce(s.study_name, '''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(samp.sample_tumor_status,'''') as `Tumor`,
    coalesce(samp.sample_type,'''') as `Analyte Type`
ORDER BY samp.sample_id '

# B4: FilesTab query - whitespace fix only (WHERE  -> WHERE)
$ws.Range("B4").Value = 'MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE f.file_type in [''XLSX'']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name,'''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id, '''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER BY f.file_name'

# Update selection to B5 (was C4)
$ws.Range("B5").Select()
